$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — update F3 and F4
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1201
$wsExhibit.Range("F4").Value = 2674

# Sheet "全部类型" (All Types) — same events are mirrored here at F5/F6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1201
$wsAll.Range("F6").Value = 2674
